# -----------------------------------------------------------------------
# Refresh the crypto price table (columns B-E) to match the latest
# GitHub Actions scrape.
#
# The sheet stores every populated cell as text (inline/shared strings) -
# even the "Price" column, whose values look numeric ("243.81", "1.001",
# "29.452.75", ...). Assigning a plain digit/dot string to Range.Value
# makes Excel auto-coerce it into a real number (losing the original
# text formatting, e.g. trailing zeros), so for any new value that is
# purely numeric-looking we temporarily force the cell to the text
# number format ("@") before writing it, then clear that formatting
# override afterwards so the cell ends up with no explicit style, same
# as before the edit.
#
# A few rows additionally had their coin name / link / price / volume
# swapped with a neighboring row (Maker<->VeChain,
# RocketPoolETH<->BabyDogeCoin, RenderToken<->EnergySwap); those are
# just more (Cell, Value) pairs below, same mechanism.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.452.75' },
    @{ Cell = 'E2'; Value = '  +0.28%  ' },
    @{ Cell = 'D3'; Value = '1.870.84' },
    @{ Cell = 'E3'; Value = '  -0.44%  ' },
    @{ Cell = 'E4'; Value = '  -0.15%  ' },
    @{ Cell = 'D5'; Value = '243.81' },
    @{ Cell = 'E5'; Value = '  +0.25%  ' },
    @{ Cell = 'D6'; Value = '0.7055' },
    @{ Cell = 'E6'; Value = '  -2.16%  ' },
    @{ Cell = 'E7'; Value = '  -0.13%  ' },
    @{ Cell = 'D8'; Value = '0.07929' },
    @{ Cell = 'E8'; Value = '  -1.37%  ' },
    @{ Cell = 'D9'; Value = '0.3139' },
    @{ Cell = 'E9'; Value = '  +0.04%  ' },
    @{ Cell = 'D10'; Value = '24.52' },
    @{ Cell = 'E10'; Value = '  -1.58%  ' },
    @{ Cell = 'D11'; Value = '0.07804' },
    @{ Cell = 'E11'; Value = '  -4.27%  ' },
    @{ Cell = 'D12'; Value = '1.867.14' },
    @{ Cell = 'D13'; Value = '93.42' },
    @{ Cell = 'E13'; Value = '  -1.22%  ' },
    @{ Cell = 'D14'; Value = '5.165' },
    @{ Cell = 'E14'; Value = '  -1.26%  ' },
    @{ Cell = 'D15'; Value = '0.7013' },
    @{ Cell = 'E15'; Value = '  -1.37%  ' },
    @{ Cell = 'D16'; Value = '6.506' },
    @{ Cell = 'E16'; Value = '  +1.18%  ' },
    @{ Cell = 'D17'; Value = '0.000008544' },
    @{ Cell = 'E17'; Value = '  +0.79%  ' },
    @{ Cell = 'D18'; Value = '29.463.49' },
    @{ Cell = 'E18'; Value = '  +0.30%  ' },
    @{ Cell = 'D19'; Value = '251.29' },
    @{ Cell = 'E19'; Value = '  +2.92%  ' },
    @{ Cell = 'D20'; Value = '2.131.93' },
    @{ Cell = 'E20'; Value = '  -0.07%  ' },
    @{ Cell = 'D21'; Value = '13.09' },
    @{ Cell = 'E21'; Value = '  -1.42%  ' },
    @{ Cell = 'D22'; Value = '0.9998' },
    @{ Cell = 'E22'; Value = '  -0.17%  ' },
    @{ Cell = 'D23'; Value = '7.603' },
    @{ Cell = 'E23'; Value = '  -1.72%  ' },
    @{ Cell = 'E24'; Value = '  -0.15%  ' },
    @{ Cell = 'D25'; Value = '0.1538' },
    @{ Cell = 'E25'; Value = '  -4.33%  ' },
    @{ Cell = 'D26'; Value = '8.999' },
    @{ Cell = 'E26'; Value = '  -0.45%  ' },
    @{ Cell = 'D27'; Value = '161.47' },
    @{ Cell = 'E27'; Value = '  -0.81%  ' },
    @{ Cell = 'D28'; Value = '18.72' },
    @{ Cell = 'E28'; Value = '  +1.19%  ' },
    @{ Cell = 'D29'; Value = '1.573' },
    @{ Cell = 'E29'; Value = '  +4.50%  ' },
    @{ Cell = 'D30'; Value = '4.304' },
    @{ Cell = 'E30'; Value = '  -2.22%  ' },
    @{ Cell = 'D31'; Value = '4.256' },
    @{ Cell = 'E31'; Value = '  -0.60%  ' },
    @{ Cell = 'E32'; Value = '  -1.76%  ' },
    @{ Cell = 'D33'; Value = '0.05268' },
    @{ Cell = 'E33'; Value = '  -1.63%  ' },
    @{ Cell = 'D34'; Value = '1.896' },
    @{ Cell = 'E34'; Value = '  -2.04%  ' },
    @{ Cell = 'D35'; Value = '0.7600' },
    @{ Cell = 'E35'; Value = '  -0.28%  ' },
    @{ Cell = 'D36'; Value = '1.183' },
    @{ Cell = 'D37'; Value = '2.706' },
    @{ Cell = 'E37'; Value = '  +0.32%  ' },
    @{ Cell = 'B38'; Value = 'VeChain' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D38'; Value = '0.01874' },
    @{ Cell = 'E38'; Value = '  +0.22%  ' },
    @{ Cell = 'B39'; Value = 'Maker' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D39'; Value = '1.278.44' },
    @{ Cell = 'E39'; Value = '  +1.33%  ' },
    @{ Cell = 'D40'; Value = '2.758' },
    @{ Cell = 'E40'; Value = '  -0.08%  ' },
    @{ Cell = 'D41'; Value = '0.8989' },
    @{ Cell = 'E41'; Value = '  -0.60%  ' },
    @{ Cell = 'D42'; Value = '109.84' },
    @{ Cell = 'E42'; Value = '  -2.91%  ' },
    @{ Cell = 'D43'; Value = '5.975' },
    @{ Cell = 'E43'; Value = '  -7.17%  ' },
    @{ Cell = 'D44'; Value = '70.81' },
    @{ Cell = 'E44'; Value = '  -4.47%  ' },
    @{ Cell = 'D45'; Value = '1.001' },
    @{ Cell = 'E45'; Value = '  -0.16%  ' },
    @{ Cell = 'B46'; Value = 'BabyDogeCoin' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' },
    @{ Cell = 'D46'; Value = '0.00000000127' },
    @{ Cell = 'E46'; Value = '  -2.92%  ' },
    @{ Cell = 'B47'; Value = 'RocketPoolETH' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' },
    @{ Cell = 'D47'; Value = '2.032.00' },
    @{ Cell = 'E47'; Value = '  +0.30%  ' },
    @{ Cell = 'B48'; Value = 'RenderToken' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D48'; Value = '1.801' },
    @{ Cell = 'E48'; Value = '  +0.17%  ' },
    @{ Cell = 'B49'; Value = 'EnergySwap' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D49'; Value = '9.615' },
    @{ Cell = 'E49'; Value = '  +1.37%  ' },
    @{ Cell = 'D50'; Value = '0.5173' },
    @{ Cell = 'E50'; Value = '  -0.44%  ' },
    @{ Cell = 'D51'; Value = '0.4300' },
    @{ Cell = 'E51'; Value = '  -0.86%  ' }
)

$numericLike = '^-?\d+(\.\d+)?$'
$textCoercedCells = @()

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Value -match $numericLike) {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $textCoercedCells += $u.Cell
    } else {
        $cell.Value = $u.Value
    }
}

foreach ($ref in $textCoercedCells) {
    $ws.Range($ref).ClearFormats()
}
